$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds a "Förändrad" (last changed) date for each data row (rows 2-426).
# Update every value in that range from 2023-09-20 (45189) to 2023-09-21 (45190).
$ws.Range("C2:C426").Value = 45190
